# Apply the Tutorial 6 attendance sheet corrections:
#  - Change the date strings in column A (rows 3-21) from DD/MM/YYYY to DD-MM-YYYY
#  - Correct a handful of attendance tally values (D/E/G/H columns) in rows 3, 5, 10, 13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A to be treated as plain text so Excel does not reinterpret
# the dash-separated strings as real dates.
$ws.Range("A3:A21").NumberFormat = "@"

$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $ws.Cells.Item($row, 1).Value = $dates[$row]
}

# Correct the attendance counts that changed alongside the date formatting.
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("H10").Value = 0

$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0
